$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update status text "Ready for handoff" -> "In Translation"
#    This string appears in:
#      - Overview sheet: E2:F3 (per-locale status columns)
#      - zh-cn sheet:     C2:C3 (Status column)
#      - de-de sheet:     C2:C3 (Status column)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# ---------------------------------------------------------------------
# 2) Shrink the (now-shorter) status columns to fit the new text:
#      - Overview sheet: columns E and F
#      - zh-cn sheet:     column C
#      - de-de sheet:     column C
#    ColumnWidth = 12.5 resolves to the display width closest to the
#    target (~13.41 chars) that Excel's column-width pixel grid allows.
# ---------------------------------------------------------------------
$wsOverview.Range("E:F").ColumnWidth = 12.5
$wsZhCn.Range("C:C").ColumnWidth = 12.5
$wsDeDe.Range("C:C").ColumnWidth = 12.5
